# The deck originally carries two theme parts:
#   ppt/theme/theme1.xml -> "Integral" / "Red Violet" color scheme (used by the slide master)
#   ppt/theme/theme2.xml -> "Office Theme" / "Office" color scheme (used by the notes master)
# The commit swaps the two themes' contents so the slide master ends up using the
# "Office Theme" palette. We reproduce that by rewriting the 12 theme colors that are
# reachable through the slide's ThemeColorScheme object (which maps onto theme1.xml,
# the theme referenced by the slide master).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Target "Office Theme" color scheme values, in ThemeColorScheme.Colors(index) order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$targetHex = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

for ($i = 0; $i -lt $targetHex.Length; $i++) {
    $hex = $targetHex[$i]
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    # PowerPoint's ColorFormat/.RGB takes a long in 0x00BBGGRR form.
    $rgbLong = ($b * 65536) + ($g * 256) + $r
    $tcs.Colors($i + 1).RGB = $rgbLong
}
